$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ALBIRENA GARCIA ANGEELO ALONSO", 146),
    @("ALAMA NIMA CLARITZA MABEL", 134),
    @("MANUEL LEUNARDO PRADO BAILON", 133),
    @("MARYURI OJEDA VALLE", 127),
    @("URRIOLA ARISMENDIZ INGRID MARYURI", 110),
    @("CORDOVA CARMEN ANGIE NATALLY", 108),
    @("AGURTO ORDINOLA LISBET JAQUELIN", 104),
    @("ROMAN GALECIO MARITZA DEL PILAR", 100),
    @("CARREÑO PALACIOS KATHERINE DE LOS MILAGROS", 99),
    @("VEGA ROBLEDO FERNANDO ERNESTO", 96),
    @("BERNAOLA CARMEN ZUMIKO YASHURY", 94),
    @("JUAREZ CARMEN PIERRE ALEXANDER", 93),
    @("ATOCHE PALACIOS LUIS ANGEL", 91),
    @("RUIDIAS FRIAS MELISSA VICTORIA", 78)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
